# CrewAI Robust Backend Ready!
# Renames header columns to the new lowercase/standardized names, realigns the
# carbon / cumulative-energy-demand figures into the correct columns, adds a
# freshly computed "climate change" column, and documents every column with a
# descriptive cell comment.

$excel.UserName = "Data Processor"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Header row: rename / re-word the column titles
# ---------------------------------------------------------------------------
$headers = @(
    @("A1", "industry"),
    @("B1", "unit"),
    @("C1", "process"),
    @("D1", "carbon (kg CO2 eq)"),
    @("E1", "ced (MJ)"),
    @("F1", "climate change (kg CO2 eq)"),
    @("G1", "region")
)

foreach ($h in $headers) {
    $ws.Range($h[0]).Value = $h[1]
}

# ---------------------------------------------------------------------------
# 2. Data rows: columns D (carbon) / E (ced) shift and F gets newly computed
#    "climate change" impact values
# ---------------------------------------------------------------------------
$data = @(
    @(2, 1.91, 0, 0.000053255868),
    @(3, 7.240696, 49.820527, 0.00020188981),
    @(4, 842.6966666666667, 11118.535, 0.023496619),
    @(5, 0.037779384, 1.4394625, 0.0000010533895),
    @(6, 0.4198036133333334, 7.2615675, 0.000011705239),
    @(7, 222.777, 3423.8066, 0.0062116139),
    @(8, 14.95335466666667, 229.09095, 0.0004169392),
    @(9, 17.26943666666667, 205.77256, 0.00048151772),
    @(10, 3.410716133333334, 37.043931, 0.000095099816),
    @(11, 0.4094384066666667, 7.1254756, 0.000011416229),
    @(12, 1.6189798, 31.632704, 0.000045141452),
    @(13, 0.32, 4.2, 0.0000089224491)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 4).Value = $row[1]
    $ws.Cells.Item($r, 5).Value = $row[2]
    $ws.Cells.Item($r, 6).Value = $row[3]
}

# ---------------------------------------------------------------------------
# 3. Header comments describing the data type held by each column
# ---------------------------------------------------------------------------
$comments = @(
    @("A1", "Data type: Categorical (text)"),
    @("B1", "Data type: Various (e.g. kg, kWh)"),
    @("C1", "Data type: Categorical (text)"),
    @("D1", "Data type: Carbon footprint"),
    @("E1", "Data type: Cumulative energy demand"),
    @("F1", "Data type: Climate change impact"),
    @("G1", "Data type: Categorical (text)")
)

foreach ($c in $comments) {
    $comment = $ws.Range($c[0]).AddComment($c[1])
    $comment.Author = "Data Processor"
}

Write-Output "applied ceramics.xlsx update"
